{"js": "const body = context.document.body;\n\n// --- Update the date paragraph (first paragraph in the document body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2023-07-05 Wednesday\", Word.InsertLocation.replace);\n\n// --- Update every multiplication-table cell, in row-major (reading) order ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n    [\"45\u00d740=\", \"64\u00d755=\", \"22\u00d720=\", \"34\u00d730=\", \"17\u00d747=\"],\n    [\"44\u00d782=\", \"37\u00d759=\", \"51\u00d711=\", \"81\u00d718=\", \"27\u00d737=\"],\n    [\"60\u00d716=\", \"28\u00d797=\", \"50\u00d736=\", \"45\u00d748=\", \"61\u00d794=\"],\n    [\"66\u00d714=\", \"69\u00d781=\", \"84\u00d745=\", \"63\u00d748=\", \"47\u00d736=\"],\n    [\"49\u00d762=\", \"81\u00d743=\", \"70\u00d771=\", \"96\u00d782=\", \"58\u00d792=\"],\n    [\"73\u00d789=\", \"37\u00d797=\", \"59\u00d721=\", \"46\u00d763=\", \"46\u00d727=\"],\n    [\"82\u00d799=\", \"80\u00d724=\", \"34\u00d728=\", \"70\u00d771=\", \"72\u00d787=\"],\n    [\"40\u00d785=\", \"15\u00d7100=\", \"96\u00d779=\", \"20\u00d775=\", \"60\u00d738=\"],\n    [\"17\u00d772=\", \"55\u00d777=\", \"59\u00d722=\", \"93\u00d762=\", \"35\u00d766=\"],\n    [\"46\u00d749=\", \"85\u00d787=\", \"25\u00d781=\", \"84\u00d735=\", \"70\u00d790=\"],\n    [\"92\u00d786=\", \"38\u00d713=\", \"38\u00d711=\", \"48\u00d714=\", \"47\u00d731=\"],\n    [\"81\u00d726=\", \"19\u00d722=\", \"44\u00d724=\", \"69\u00d729=\", \"53\u00d714=\"],\n    [\"68\u00d797=\", \"74\u00d752=\", \"95\u00d754=\", \"12\u00d738=\", \"89\u00d789=\"],\n    [\"24\u00d770=\", \"68\u00d762=\", \"41\u00d778=\", \"37\u00d738=\", \"100\u00d713=\"],\n    [\"33\u00d752=\", \"97\u00d730=\", \"27\u00d717=\", \"49\u00d713=\", \"88\u00d791=\"],\n    [\"52\u00d760=\", \"36\u00d751=\", \"67\u00d715=\", \"99\u00d742=\", \"99\u00d767=\"],\n    [\"49\u00d725=\", \"21\u00d779=\", \"47\u00d718=\", \"60\u00d737=\", \"11\u00d768=\"],\n    [\"88\u00d799=\", \"29\u00d751=\", \"36\u00d7100=\", \"96\u00d790=\", \"41\u00d787=\"],\n    [\"95\u00d711=\", \"16\u00d743=\", \"72\u00d720=\", \"38\u00d716=\", \"77\u00d742=\"],\n    [\"62\u00d726=\", \"62\u00d756=\", \"100\u00d755=\", \"57\u00d759=\", \"74\u00d718=\"]\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Update the date heading paragraph (first paragraph in the body) ---\n$dateParagraph = $d.Paragraphs.Item(1)\n$dateParagraph.Range.Text = \"2023-07-05 Wednesday\"\n\n# --- Update every multiplication-table cell, in row-major order ---\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"45\u00d740=\", \"64\u00d755=\", \"22\u00d720=\", \"34\u00d730=\", \"17\u00d747=\"),\n    @(\"44\u00d782=\", \"37\u00d759=\", \"51\u00d711=\", \"81\u00d718=\", \"27\u00d737=\"),\n    @(\"60\u00d716=\", \"28\u00d797=\", \"50\u00d736=\", \"45\u00d748=\", \"61\u00d794=\"),\n    @(\"66\u00d714=\", \"69\u00d781=\", \"84\u00d745=\", \"63\u00d748=\", \"47\u00d736=\"),\n    @(\"49\u00d762=\", \"81\u00d743=\", \"70\u00d771=\", \"96\u00d782=\", \"58\u00d792=\"),\n    @(\"73\u00d789=\", \"37\u00d797=\", \"59\u00d721=\", \"46\u00d763=\", \"46\u00d727=\"),\n    @(\"82\u00d799=\", \"80\u00d724=\", \"34\u00d728=\", \"70\u00d771=\", \"72\u00d787=\"),\n    @(\"40\u00d785=\", \"15\u00d7100=\", \"96\u00d779=\", \"20\u00d775=\", \"60\u00d738=\"),\n    @(\"17\u00d772=\", \"55\u00d777=\", \"59\u00d722=\", \"93\u00d762=\", \"35\u00d766=\"),\n    @(\"46\u00d749=\", \"85\u00d787=\", \"25\u00d781=\", \"84\u00d735=\", \"70\u00d790=\"),\n    @(\"92\u00d786=\", \"38\u00d713=\", \"38\u00d711=\", \"48\u00d714=\", \"47\u00d731=\"),\n    @(\"81\u00d726=\", \"19\u00d722=\", \"44\u00d724=\", \"69\u00d729=\", \"53\u00d714=\"),\n    @(\"68\u00d797=\", \"74\u00d752=\", \"95\u00d754=\", \"12\u00d738=\", \"89\u00d789=\"),\n    @(\"24\u00d770=\", \"68\u00d762=\", \"41\u00d778=\", \"37\u00d738=\", \"100\u00d713=\"),\n    @(\"33\u00d752=\", \"97\u00d730=\", \"27\u00d717=\", \"49\u00d713=\", \"88\u00d791=\"),\n    @(\"52\u00d760=\", \"36\u00d751=\", \"67\u00d715=\", \"99\u00d742=\", \"99\u00d767=\"),\n    @(\"49\u00d725=\", \"21\u00d779=\", \"47\u00d718=\", \"60\u00d737=\", \"11\u00d768=\"),\n    @(\"88\u00d799=\", \"29\u00d751=\", \"36\u00d7100=\", \"96\u00d790=\", \"41\u00d787=\"),\n    @(\"95\u00d711=\", \"16\u00d743=\", \"72\u00d720=\", \"38\u00d716=\", \"77\u00d742=\"),\n    @(\"62\u00d726=\", \"62\u00d756=\", \"100\u00d755=\", \"57\u00d759=\", \"74\u00d718=\")\n)\n\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n\n"}
